$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

# Rows 4-6: cyclic relabel/reorder of "Extreme Low Flow Method" entries (A-F)
# and their associated Flow Location / Flow Type / Min / Max / Strategy columns.
# New row 4 <= old row 5 data (relabelled C.)
$ws.Range("A4").Value = "C. 85%, 65%, and 50% of  2000 to 2018 average flow"
$ws.Range("B4").Value = "Lake Powell"
$ws.Range("C4").Value = "Regulated Inflow"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Release 95% of regulated  inflow."

# New row 5 <= old row 6 data (relabelled D.)
$ws.Range("A5").Value = "D. Reclamation's Post 2026  web tool"
$ws.Range("B5").Value = "Lake Powell"
$ws.Range("C5").Value = "Release"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = "Release to prevent drawdown  to 3,490 feet."

# New row 6 <= old row 4 data (relabelled E.)
$ws.Range("A6").Value = "E. Low Lake Powell releases + gains  through Grand Canyon"
$ws.Range("B6").Value = "Lake Mead"
$ws.Range("C6").Value = "Regulated Inflow"
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "Rule curve; Consumption equals or less  than inflow minus evaporation."

# Row 7 label changed
$ws.Range("A7").Value = '$Flo'

# Update the sheet's selection to A7
$ws.Range("A7").Select()
